# Update Diebold-Mariano statistic (column C) and P-Value (column D)
# for rows 2-11 per corrected values from the commit
# "Correcion a Diebold Mariano y revision de Cap1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 2;  C = -0.07983774705404588; D = 0.9370881031581431 },
    @{ Row = 3;  C = -0.580263322589555;   D = 0.5676334687541575 },
    @{ Row = 4;  C = -2.184598312191488;   D = 0.03985905858210481 },
    @{ Row = 5;  C = -0.9126187327275742;  D = 0.3713361602658627 },
    @{ Row = 6;  C = -0.5672919854961624;  D = 0.5762593634827216 },
    @{ Row = 7;  C = -1.54701944760919;    D = 0.1361240161137447 },
    @{ Row = 8;  C = -1.051001375472754;   D = 0.3046662384684842 },
    @{ Row = 9;  C = -1.234111707097367;   D = 0.2301825960583397 },
    @{ Row = 10; C = -0.8000733981514283;  D = 0.4322232391561178 },
    @{ Row = 11; C = 0.4692201312574588;   D = 0.643529163473457 }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
